# Apply updated crypto price/volume values (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.114.85"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "3.529.51"
$ws.Range("E3").Value = "  +3.37%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'596.44"
$ws.Range("E5").Value = "  +2.15%  "
$ws.Range("D6").Value = "'137.51"
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("D7").Value = "3.530.08"
$ws.Range("E7").Value = "  +3.38%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("D9").Value = "'0.492"
$ws.Range("E9").Value = "  -0.77%  "
$ws.Range("E10").Value = "  +2.68%  "
$ws.Range("D11").Value = "'6.90"
$ws.Range("E11").Value = "  -4.70%  "
$ws.Range("E12").Value = "  +2.72%  "
$ws.Range("D13").Value = "4.130.17"
$ws.Range("E13").Value = "  +3.66%  "
$ws.Range("E14").Value = "  +1.83%  "
$ws.Range("D15").Value = "'27.06"
$ws.Range("E15").Value = "  +4.02%  "
$ws.Range("D16").Value = "3.532.35"
$ws.Range("E16").Value = "  +3.36%  "
$ws.Range("E17").Value = "  +1.20%  "
$ws.Range("D18").Value = "65.109.47"
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("D19").Value = "'10.11"
$ws.Range("E19").Value = "  +3.81%  "
$ws.Range("D20").Value = "'5.89"
$ws.Range("E20").Value = "  +0.54%  "
$ws.Range("D21").Value = "'14.22"
$ws.Range("E21").Value = "  +5.01%  "
$ws.Range("D22").Value = "'390.02"
$ws.Range("E22").Value = "  +1.34%  "
$ws.Range("E23").Value = "  +3.01%  "
$ws.Range("D24").Value = "3.671.44"
$ws.Range("E24").Value = "  +3.49%  "
$ws.Range("D25").Value = "'73.38"
$ws.Range("E25").Value = "  +1.28%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").Value = "'0.0000113"
$ws.Range("E27").Value = "  +7.89%  "
$ws.Range("D28").Value = "'7.77"
$ws.Range("E28").Value = "  +10.28%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("E30").Value = "  +2.96%  "
$ws.Range("E31").Value = "  +0.79%  "
$ws.Range("D32").Value = "3.551.68"
$ws.Range("E32").Value = "  +3.98%  "
$ws.Range("D34").Value = "'23.75"
$ws.Range("E34").Value = "  +4.44%  "
$ws.Range("E35").Value = "  +0.77%  "
$ws.Range("E36").Value = "  +14.47%  "
$ws.Range("D37").Value = "'169.86"
$ws.Range("E37").Value = "  -0.39%  "
$ws.Range("E38").Value = "  +7.31%  "
$ws.Range("E39").Value = "  +2.13%  "
$ws.Range("E40").Value = "  +5.27%  "
$ws.Range("E41").Value = "  +5.39%  "
$ws.Range("D42").Value = "'0.823"
$ws.Range("E42").Value = "  +1.56%  "
$ws.Range("D43").Value = "'26.28"
$ws.Range("E43").Value = "  +17.78%  "
$ws.Range("D44").Value = "'42.60"
$ws.Range("E44").Value = "  -2.34%  "
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("E46").Value = "  +0.68%  "
$ws.Range("E47").Value = "  +8.66%  "
$ws.Range("D48").Value = "'1.66"
$ws.Range("E48").Value = "  +4.04%  "
$ws.Range("D49").Value = "'6.77"
$ws.Range("E49").Value = "  +3.97%  "
$ws.Range("D50").Value = "2.382.52"
$ws.Range("E50").Value = "  +9.70%  "
$ws.Range("E51").Value = "  +19.25%  "
